$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8 data (20 minute trade where close price could not be located from Yahoo)
$ws.Range("A8").Value = 9580.3799999999992
$ws.Range("B8").Value = 9697.7199999999993
$ws.Range("C8").Value = 309.02999999999997
$ws.Range("D8").Value = 305.29000000000002
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = -1.21
$ws.Range("G8").Value = 42608.617951388886
$ws.Range("G8").NumberFormat = "m/d/yy h:mm"
$ws.Range("H8").Value = $false
